$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.017531666666667
$ws.Range("H2").Value = 27.052595
$ws.Range("I2").Value = 0.09025733169883475
$ws.Range("J2").Value = 0.09025733169883476
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 26.93702966666667
$ws.Range("N2").Value = 80.81108900000001
$ws.Range("O2").Value = 0.05621802850548585
$ws.Range("P2").Value = 0.05621802850548584
$ws.Range("Q2").Value = 242.9055180251061
$ws.Range("R2").Value = 2186.149662225955
$ws.Range("S2").Value = 0.005074089246274184
$ws.Range("T2").Value = 0.005074089246274184

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.017531666666667
$ws.Range("H3").Value = 27.052595
$ws.Range("I3").Value = 0.09025733169883475
$ws.Range("J3").Value = 0.09025733169883476
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 191.7798513333333
$ws.Range("N3").Value = 575.339554
$ws.Range("O3").Value = 0.4002477363856031
$ws.Range("P3").Value = 0.4002477363856031
$ws.Range("Q3").Value = 1729.380882426959
$ws.Range("R3").Value = 15564.42794184263
$ws.Range("S3").Value = 0.03612529270466315
$ws.Range("T3").Value = 0.03612529270466315

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.017531666666667
$ws.Range("H4").Value = 27.052595
$ws.Range("I4").Value = 0.09025733169883475
$ws.Range("J4").Value = 0.09025733169883476
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 172.1284333333333
$ws.Range("N4").Value = 516.3853
$ws.Range("O4").Value = 0.3592349004876528
$ws.Range("P4").Value = 0.3592349004876528
$ws.Range("Q4").Value = 1552.173598317056
$ws.Range("R4").Value = 13969.5623848535
$ws.Range("S4").Value = 0.03242358357111198
$ws.Range("T4").Value = 0.03242358357111198

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.017531666666667
$ws.Range("H5").Value = 27.052595
$ws.Range("I5").Value = 0.09025733169883475
$ws.Range("J5").Value = 0.09025733169883476
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 88.307555
$ws.Range("N5").Value = 264.922665
$ws.Range("O5").Value = 0.1842993346212582
$ws.Range("P5").Value = 0.1842993346212582
$ws.Range("Q5").Value = 796.3161736184082
$ws.Range("R5").Value = 7166.845562565675
$ws.Range("S5").Value = 0.01663436617678544
$ws.Range("T5").Value = 0.01663436617678544

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 73.986959
$ws.Range("H6").Value = 221.960877
$ws.Range("I6").Value = 0.7405425061645015
$ws.Range("J6").Value = 0.7405425061645016
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 26.93702966666667
$ws.Range("N6").Value = 80.81108900000001
$ws.Range("O6").Value = 0.05621802850548585
$ws.Range("P6").Value = 0.05621802850548584
$ws.Range("Q6").Value = 1992.988909529451
$ws.Range("R6").Value = 17936.90018576505
$ws.Range("S6").Value = 0.04163183972107987
$ws.Range("T6").Value = 0.04163183972107987

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 73.986959
$ws.Range("H7").Value = 221.960877
$ws.Range("I7").Value = 0.7405425061645015
$ws.Range("J7").Value = 0.7405425061645016
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 191.7798513333333
$ws.Range("N7").Value = 575.339554
$ws.Range("O7").Value = 0.4002477363856031
$ws.Range("P7").Value = 0.4002477363856031
$ws.Range("Q7").Value = 14189.20799762543
$ws.Range("R7").Value = 127702.8719786289
$ws.Range("S7").Value = 0.2964004617896633
$ws.Range("T7").Value = 0.2964004617896633

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 73.986959
$ws.Range("H8").Value = 221.960877
$ws.Range("I8").Value = 0.7405425061645015
$ws.Range("J8").Value = 0.7405425061645016
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 172.1284333333333
$ws.Range("N8").Value = 516.3853
$ws.Range("O8").Value = 0.3592349004876528
$ws.Range("P8").Value = 0.3592349004876528
$ws.Range("Q8").Value = 12735.25933976757
$ws.Range("R8").Value = 114617.3340579081
$ws.Range("S8").Value = 0.2660287135088817
$ws.Range("T8").Value = 0.2660287135088817

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 73.986959
$ws.Range("H9").Value = 221.960877
$ws.Range("I9").Value = 0.7405425061645015
$ws.Range("J9").Value = 0.7405425061645016
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 88.307555
$ws.Range("N9").Value = 264.922665
$ws.Range("O9").Value = 0.1842993346212582
$ws.Range("P9").Value = 0.1842993346212582
$ws.Range("Q9").Value = 6533.607451175244
$ws.Range("R9").Value = 58802.4670605772
$ws.Range("S9").Value = 0.1364814911448766
$ws.Range("T9").Value = 0.1364814911448766

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.075612333333333
$ws.Range("H10").Value = 3.226837
$ws.Range("I10").Value = 0.01076590609688545
$ws.Range("J10").Value = 0.01076590609688545
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 26.93702966666667
$ws.Range("N10").Value = 80.81108900000001
$ws.Range("O10").Value = 0.05621802850548585
$ws.Range("P10").Value = 0.05621802850548584
$ws.Range("Q10").Value = 28.97380133283255
$ws.Range("R10").Value = 260.764211995493
$ws.Range("S10").Value = 0.00060523801584209
$ws.Range("T10").Value = 0.00060523801584209

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.075612333333333
$ws.Range("H11").Value = 3.226837
$ws.Range("I11").Value = 0.01076590609688545
$ws.Range("J11").Value = 0.01076590609688545
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 191.7798513333333
$ws.Range("N11").Value = 575.339554
$ws.Range("O11").Value = 0.4002477363856031
$ws.Range("P11").Value = 0.4002477363856031
$ws.Range("Q11").Value = 206.2807733789664
$ws.Range("R11").Value = 1856.526960410698
$ws.Range("S11").Value = 0.004309029545418364
$ws.Range("T11").Value = 0.004309029545418364

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.075612333333333
$ws.Range("H12").Value = 3.226837
$ws.Range("I12").Value = 0.01076590609688545
$ws.Range("J12").Value = 0.01076590609688545
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 172.1284333333333
$ws.Range("N12").Value = 516.3853
$ws.Range("O12").Value = 0.3592349004876528
$ws.Range("P12").Value = 0.3592349004876528
$ws.Range("Q12").Value = 185.1434658106778
$ws.Range("R12").Value = 1666.2911922961
$ws.Range("S12").Value = 0.003867489205374059
$ws.Range("T12").Value = 0.003867489205374059

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.075612333333333
$ws.Range("H13").Value = 3.226837
$ws.Range("I13").Value = 0.01076590609688545
$ws.Range("J13").Value = 0.01076590609688545
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 88.307555
$ws.Range("N13").Value = 264.922665
$ws.Range("O13").Value = 0.1842993346212582
$ws.Range("P13").Value = 0.1842993346212582
$ws.Range("Q13").Value = 94.98469528451164
$ws.Range("R13").Value = 854.8622575606049
$ws.Range("S13").Value = 0.001984149330250935
$ws.Range("T13").Value = 0.001984149330250935

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.829029
$ws.Range("H14").Value = 47.487087
$ws.Range("I14").Value = 0.1584342560397782
$ws.Range("J14").Value = 0.1584342560397782
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 26.93702966666667
$ws.Range("N14").Value = 80.81108900000001
$ws.Range("O14").Value = 0.05621802850548585
$ws.Range("P14").Value = 0.05621802850548584
$ws.Range("Q14").Value = 426.3870237675271
$ws.Range("R14").Value = 3837.483213907743
$ws.Range("S14").Value = 0.008906861522289694
$ws.Range("T14").Value = 0.008906861522289696

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.829029
$ws.Range("H15").Value = 47.487087
$ws.Range("I15").Value = 0.1584342560397782
$ws.Range("J15").Value = 0.1584342560397782
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 191.7798513333333
$ws.Range("N15").Value = 575.339554
$ws.Range("O15").Value = 0.4002477363856031
$ws.Range("P15").Value = 0.4002477363856031
$ws.Range("Q15").Value = 3035.688828371022
$ws.Range("R15").Value = 27321.1994553392
$ws.Range("S15").Value = 0.0634129523458583
$ws.Range("T15").Value = 0.0634129523458583

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.829029
$ws.Range("H16").Value = 47.487087
$ws.Range("I16").Value = 0.1584342560397782
$ws.Range("J16").Value = 0.1584342560397782
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 172.1284333333333
$ws.Range("N16").Value = 516.3853
$ws.Range("O16").Value = 0.3592349004876528
$ws.Range("P16").Value = 0.3592349004876528
$ws.Range("Q16").Value = 2724.6259629579
$ws.Range("R16").Value = 24521.6336666211
$ws.Range("S16").Value = 0.05691511420228503
$ws.Range("T16").Value = 0.05691511420228505

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.829029
$ws.Range("H17").Value = 47.487087
$ws.Range("I17").Value = 0.1584342560397782
$ws.Range("J17").Value = 0.1584342560397782
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 88.307555
$ws.Range("N17").Value = 264.922665
$ws.Range("O17").Value = 0.1842993346212582
$ws.Range("P17").Value = 0.1842993346212582
$ws.Range("Q17").Value = 1397.822849014095
$ws.Range("R17").Value = 12580.40564112686
$ws.Range("S17").Value = 0.02919932796934517
$ws.Range("T17").Value = 0.02919932796934518
